$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.491435050964355
$ws.Range("B1").Value = 2.326718330383301
$ws.Range("C1").Value = 5.221031188964844
$ws.Range("D1").Value = 3.45506763458252
$ws.Range("E1").Value = 1.056059598922729
